$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H8").Value = 30.785715
$ws.Range("I8").Value = 33.416668
$ws.Range("K8").Value = 100.250004
$ws.Range("M8").Value = 38.749996
$ws.Range("H17").Value = 1549.8334
$ws.Range("J17").Value = 1549.8334
$ws.Range("L17").Value = 4649.5002
$ws.Range("N17").Value = -4985.5002
$ws.Range("H25").Value = 1080
$ws.Range("J25").Value = 1080
$ws.Range("L25").Value = 3240
$ws.Range("N25").Value = -3588
$ws.Range("H28").Value = 202.66667
$ws.Range("I28").Value = 202.66667
$ws.Range("K28").Value = 202.66667
$ws.Range("M28").Value = 282.33333
$ws.Range("H41").Value = 267.125
$ws.Range("I41").Value = 233.6
$ws.Range("K41").Value = 233.6
$ws.Range("M41").Value = 206.4
$ws.Range("H111").Value = 409.66666
$ws.Range("I111").Value = 409.66666
$ws.Range("K111").Value = 1228.99998
$ws.Range("M111").Value = 1838.00002
$ws.Range("H115").Value = 380
$ws.Range("I115").Value = 380
$ws.Range("K115").Value = 1140
$ws.Range("M115").Value = 427
$ws.Range("H116").Value = 6348.75
$ws.Range("I116").Value = 5959.4
$ws.Range("J116").Value = 6997.6665
$ws.Range("K116").Value = 5959.4
$ws.Range("L116").Value = 6997.6665
$ws.Range("M116").Value = -2517.4
$ws.Range("N116").Value = -13881.6665
$ws.Range("H137").Value = 1629.6897
$ws.Range("I137").Value = 1369.72
$ws.Range("K137").Value = 4109.16
$ws.Range("M137").Value = -1559.16
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1170
$ws.Range("I2").Value = 1115
$ws.Range("K2").Value = 1115
$ws.Range("M2").Value = -1002
$ws.Range("H32").Value = 1831.6613
$ws.Range("I32").Value = 1687.5254
$ws.Range("K32").Value = 1687.5254
$ws.Range("M32").Value = -1400.5254
$ws.Range("H61").Value = 2421.5
$ws.Range("I61").Value = 2562.7778
$ws.Range("J61").Value = 1150
$ws.Range("K61").Value = 2562.7778
$ws.Range("L61").Value = 1150
$ws.Range("M61").Value = -2350.7778
$ws.Range("N61").Value = -1574
$ws.Range("H74").Value = 853.1429000000001
$ws.Range("I74").Value = 662.4167
$ws.Range("K74").Value = 662.4167
$ws.Range("M74").Value = 211.5833
$ws.Range("H77").Value = 853.1429000000001
$ws.Range("I77").Value = 662.4167
$ws.Range("K77").Value = 3312.0835
$ws.Range("M77").Value = 1055.9165
$ws.Range("H88").Value = 2273.4167
$ws.Range("I88").Value = 1189
$ws.Range("K88").Value = 1189
$ws.Range("M88").Value = -783
$ws.Range("H91").Value = 2273.4167
$ws.Range("I91").Value = 1189
$ws.Range("K91").Value = 1189
$ws.Range("M91").Value = 215
$ws.Range("H116").Value = 1170
$ws.Range("I116").Value = 1115
$ws.Range("K116").Value = 1115
$ws.Range("M116").Value = 1179
$ws.Range("H136").Value = 2421.5
$ws.Range("I136").Value = 2562.7778
$ws.Range("J136").Value = 1150
$ws.Range("K136").Value = 7688.3334
$ws.Range("L136").Value = 3450
$ws.Range("M136").Value = -5138.3334
$ws.Range("N136").Value = -8550
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1170
$ws.Range("I3").Value = 1115
$ws.Range("K3").Value = 1115
$ws.Range("M3").Value = -1001
$ws.Range("H86").Value = 3703.4
$ws.Range("I86").Value = 2508.6
$ws.Range("J86").Value = 4898.2
$ws.Range("K86").Value = 2508.6
$ws.Range("L86").Value = 4898.2
$ws.Range("M86").Value = -1385.6
$ws.Range("N86").Value = -7144.2
$ws.Range("H89").Value = 3703.4
$ws.Range("I89").Value = 2508.6
$ws.Range("J89").Value = 4898.2
$ws.Range("K89").Value = 12543
$ws.Range("L89").Value = 24491
$ws.Range("M89").Value = -6927
$ws.Range("N89").Value = -35723
$ws.Range("H134").Value = 10913.833
$ws.Range("I134").Value = 10896.7
$ws.Range("K134").Value = 32690.1
$ws.Range("M134").Value = -30155.1
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 3333.9285
$ws.Range("I58").Value = 2717.5
$ws.Range("J58").Value = 4875
$ws.Range("K58").Value = 2717.5
$ws.Range("L58").Value = 4875
$ws.Range("M58").Value = -2514.5
$ws.Range("N58").Value = -5281
$ws.Range("H132").Value = 2932.6667
$ws.Range("I132").Value = 2932.6667
$ws.Range("K132").Value = 8798.000100000001
$ws.Range("M132").Value = -6268.000100000001
$ws.Range("H136").Value = 3333.9285
$ws.Range("I136").Value = 2717.5
$ws.Range("J136").Value = 4875
$ws.Range("K136").Value = 8152.5
$ws.Range("L136").Value = 14625
$ws.Range("M136").Value = -5602.5
$ws.Range("N136").Value = -19725
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H33").Value = 382.33334
$ws.Range("I33").Value = 150
$ws.Range("J33").Value = 498.5
$ws.Range("K33").Value = 900
$ws.Range("L33").Value = 2991
$ws.Range("M33").Value = -617
$ws.Range("N33").Value = -3557
$ws.Range("H113").Value = 924.82355
$ws.Range("I113").Value = 580.2222
$ws.Range("J113").Value = 1312.5
$ws.Range("K113").Value = 1740.6666
$ws.Range("L113").Value = 3937.5
$ws.Range("M113").Value = 429.3334
$ws.Range("N113").Value = -8277.5
$ws.Range("H122").Value = 1262.8334
$ws.Range("J122").Value = 1545.75
$ws.Range("L122").Value = 13911.75
$ws.Range("N122").Value = -18811.75
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H22").Value = 44.5
$ws.Range("I22").Value = 0
$ws.Range("K22").Value = 0
$ws.Range("M22").Value = ""
$ws.Range("H113").Value = 3430
$ws.Range("I113").Value = 300
$ws.Range("K113").Value = 300
$ws.Range("M113").Value = 1870
$ws.Range("H132").Value = 2274.6924
$ws.Range("I132").Value = 2273.84
$ws.Range("K132").Value = 6821.52
$ws.Range("M132").Value = -4291.52
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H4").Value = 0
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 0
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = ""
$ws.Range("N4").Value = ""
$ws.Range("H16").Value = 860.6
$ws.Range("I16").Value = 701
$ws.Range("K16").Value = 701
$ws.Range("M16").Value = -531
$ws.Range("H22").Value = 1797
$ws.Range("J22").Value = 1995.25
$ws.Range("L22").Value = 1995.25
$ws.Range("N22").Value = -2585.25
$ws.Range("H24").Value = 50007
$ws.Range("J24").Value = 50007
$ws.Range("L24").Value = 50007
$ws.Range("N24").Value = -50693
$ws.Range("H27").Value = 1797
$ws.Range("J27").Value = 1995.25
$ws.Range("L27").Value = 1995.25
$ws.Range("N27").Value = -2209.25
$ws.Range("H28").Value = 0
$ws.Range("I28").Value = 0
$ws.Range("J28").Value = 0
$ws.Range("K28").Value = 0
$ws.Range("L28").Value = 0
$ws.Range("M28").Value = ""
$ws.Range("N28").Value = ""
$ws.Range("H31").Value = 1412.75
$ws.Range("I31").Value = 1475.5
$ws.Range("K31").Value = 1475.5
$ws.Range("M31").Value = -1227.5
$ws.Range("H37").Value = 0
$ws.Range("I37").Value = 0
$ws.Range("J37").Value = 0
$ws.Range("K37").Value = 0
$ws.Range("L37").Value = 0
$ws.Range("M37").Value = ""
$ws.Range("N37").Value = ""
$ws.Range("H39").Value = 20810.4
$ws.Range("I39").Value = 500
$ws.Range("K39").Value = 500
$ws.Range("M39").Value = -40
$ws.Range("H40").Value = 3000
$ws.Range("I40").Value = 3000
$ws.Range("J40").Value = 0
$ws.Range("K40").Value = 3000
$ws.Range("L40").Value = 0
$ws.Range("M40").Value = -2864
$ws.Range("N40").Value = ""
$ws.Range("H55").Value = 192.33333
$ws.Range("J55").Value = 192.75
$ws.Range("L55").Value = 192.75
$ws.Range("N55").Value = -538.75
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H18").Value = 50007
$ws.Range("I18").Value = 0
$ws.Range("J18").Value = 50007
$ws.Range("K18").Value = 0
$ws.Range("L18").Value = 50007
$ws.Range("M18").Value = ""
$ws.Range("N18").Value = -50353
$ws.Range("H100").Value = 1159.7142
$ws.Range("I100").Value = 1261.6666
$ws.Range("K100").Value = 2523.3332
$ws.Range("M100").Value = -1982.3332
